{"js": "// Replace the date line and the two-digit multiplication problems with\n// the new values from the commit, matching each old string exactly once.\nconst replacements = [\n  [\"2025-08-31 Sunday\", \"2025-09-01 Monday\"],\n  [\"30\u00d752=\", \"85\u00d744=\"],\n  [\"76\u00d782=\", \"39\u00d799=\"],\n  [\"27\u00d784=\", \"90\u00d750=\"],\n  [\"63\u00d727=\", \"58\u00d723=\"],\n  [\"21\u00d754=\", \"62\u00d731=\"],\n  [\"39\u00d744=\", \"26\u00d723=\"],\n  [\"16\u00d786=\", \"44\u00d798=\"],\n  [\"98\u00d772=\", \"71\u00d746=\"],\n  [\"33\u00d714=\", \"27\u00d774=\"],\n  [\"61\u00d730=\", \"73\u00d732=\"],\n  [\"38\u00d738=\", \"59\u00d728=\"],\n  [\"67\u00d771=\", \"89\u00d772=\"],\n  [\"56\u00d722=\", \"92\u00d798=\"],\n  [\"77\u00d780=\", \"73\u00d796=\"],\n  [\"89\u00d760=\", \"39\u00d730=\"],\n  [\"27\u00d721=\", \"64\u00d739=\"],\n  [\"80\u00d768=\", \"65\u00d779=\"],\n  [\"67\u00d770=\", \"88\u00d720=\"],\n  [\"37\u00d727=\", \"52\u00d757=\"],\n  [\"42\u00d775=\", \"81\u00d723=\"],\n  [\"75\u00d754=\", \"98\u00d728=\"],\n  [\"19\u00d765=\", \"29\u00d796=\"],\n  [\"65\u00d717=\", \"30\u00d732=\"],\n  [\"39\u00d716=\", \"39\u00d792=\"],\n  [\"85\u00d735=\", \"53\u00d775=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-08-31 Sunday\"; New = \"2025-09-01 Monday\" },\n    @{ Old = \"30\u00d752=\"; New = \"85\u00d744=\" },\n    @{ Old = \"76\u00d782=\"; New = \"39\u00d799=\" },\n    @{ Old = \"27\u00d784=\"; New = \"90\u00d750=\" },\n    @{ Old = \"63\u00d727=\"; New = \"58\u00d723=\" },\n    @{ Old = \"21\u00d754=\"; New = \"62\u00d731=\" },\n    @{ Old = \"39\u00d744=\"; New = \"26\u00d723=\" },\n    @{ Old = \"16\u00d786=\"; New = \"44\u00d798=\" },\n    @{ Old = \"98\u00d772=\"; New = \"71\u00d746=\" },\n    @{ Old = \"33\u00d714=\"; New = \"27\u00d774=\" },\n    @{ Old = \"61\u00d730=\"; New = \"73\u00d732=\" },\n    @{ Old = \"38\u00d738=\"; New = \"59\u00d728=\" },\n    @{ Old = \"67\u00d771=\"; New = \"89\u00d772=\" },\n    @{ Old = \"56\u00d722=\"; New = \"92\u00d798=\" },\n    @{ Old = \"77\u00d780=\"; New = \"73\u00d796=\" },\n    @{ Old = \"89\u00d760=\"; New = \"39\u00d730=\" },\n    @{ Old = \"27\u00d721=\"; New = \"64\u00d739=\" },\n    @{ Old = \"80\u00d768=\"; New = \"65\u00d779=\" },\n    @{ Old = \"67\u00d770=\"; New = \"88\u00d720=\" },\n    @{ Old = \"37\u00d727=\"; New = \"52\u00d757=\" },\n    @{ Old = \"42\u00d775=\"; New = \"81\u00d723=\" },\n    @{ Old = \"75\u00d754=\"; New = \"98\u00d728=\" },\n    @{ Old = \"19\u00d765=\"; New = \"29\u00d796=\" },\n    @{ Old = \"65\u00d717=\"; New = \"30\u00d732=\" },\n    @{ Old = \"39\u00d716=\"; New = \"39\u00d792=\" },\n    @{ Old = \"85\u00d735=\"; New = \"53\u00d775=\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
